# Auto-generated edit script applying scheduled market-data refresh values
# to the currentAveragePrice / Leve price / profit columns (H:N) across all
# 8 job sheets, per the authoritative diff.

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ALC_values = @{
  "H4" = 237
  "I4" = 257.66666
  "K4" = 257.66666
  "M4" = -143.66666
  "H28" = 2084.9048
  "I28" = 1988.8
  "K28" = 1988.8
  "M28" = -1503.8
  "H33" = 373.5
  "I33" = 408.57144
  "K33" = 408.57144
  "M33" = -179.57144
  "H62" = 3531.3914
  "I62" = 2966.7334
  "J62" = 4590.125
  "K62" = 2966.7334
  "L62" = 4590.125
  "M62" = -2342.7334
  "N62" = -5838.125
  "H65" = 3531.3914
  "I65" = 2966.7334
  "J65" = 4590.125
  "K65" = 14833.667
  "L65" = 22950.625
  "M65" = -11713.667
  "N65" = -29190.625
  "H69" = 55561556
  "I69" = 166669170
  "J69" = 7750
  "K69" = 500007510
  "L69" = 23250
  "M69" = -500006636
  "N69" = -24998
  "H72" = 55561556
  "I72" = 166669170
  "J72" = 7750
  "K72" = 1500022530
  "L72" = 69750
  "M72" = -1500018162
  "N72" = -78486
  "H107" = 890.9
  "J107" = 98
  "L107" = 98
  "N107" = -3938
  "H113" = 1200
  "J113" = 1200
  "L113" = 1200
  "N113" = -7708
  "H129" = 2086.325
  "J129" = 2396.303
  "L129" = 7188.909
  "N129" = -17188.909
  "H132" = 4764.2856
  "I132" = 1324.8864
  "J132" = 35031
  "K132" = 3974.6592
  "L132" = 105093
  "M132" = -1444.6592
  "N132" = -110153
  "H135" = 3092.4443
  "I135" = 944.6667
  "K135" = 8502.0003
  "M135" = -5967.0003
  "H137" = 6360578
  "I137" = 305335.94
  "J137" = 15875958
  "K137" = 916007.8200000001
  "L137" = 47627874
  "M137" = -913457.8200000001
  "N137" = -47632974
  "H138" = 6592.5
  "I138" = 2423.8333
  "J138" = 7250.7104
  "K138" = 7271.499899999999
  "L138" = 21752.1312
  "M138" = -2131.499899999999
  "N138" = -32032.1312
  "H141" = 2142.4
  "I141" = 2142.4
  "J141" = 0
  "K141" = 6427.200000000001
  "L141" = 0
  "M141" = -1247.200000000001
}
foreach ($ref in $ALC_values.Keys) {
    $ws.Range($ref).Value = $ALC_values[$ref]
}
$ALC_deletes = @("N141")
foreach ($ref in $ALC_deletes) {
    $ws.Range($ref).ClearContents()
}

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ARM_values = @{
  "H2" = 898.84
  "I2" = 821.5454999999999
  "K2" = 821.5454999999999
  "M2" = -708.5454999999999
  "H32" = 15691.271
  "I32" = 15403.048
  "K32" = 15403.048
  "M32" = -15116.048
  "H45" = 65064.723
  "I45" = 88435.84
  "J45" = 4299.8
  "K45" = 88435.84
  "L45" = 4299.8
  "M45" = -88058.84
  "N45" = -5053.8
  "H61" = 4073.6667
  "I61" = 3832.1458
  "K61" = 3832.1458
  "M61" = -3620.1458
  "H74" = 957.34485
  "I74" = 922.25
  "J74" = 1940
  "K74" = 922.25
  "L74" = 1940
  "M74" = -48.25
  "N74" = -3688
  "H77" = 957.34485
  "I77" = 922.25
  "J77" = 1940
  "K77" = 4611.25
  "L77" = 9700
  "M77" = -243.25
  "N77" = -18436
  "H109" = 59500
  "J109" = 59500
  "L109" = 59500
  "N109" = -62274
  "H116" = 898.84
  "I116" = 821.5454999999999
  "K116" = 821.5454999999999
  "M116" = 1472.4545
  "H124" = 20333
  "J124" = 20333
  "L124" = 20333
  "N124" = -30153
  "H132" = 18952.15
  "I132" = 23819.23
  "J132" = 9913.286
  "K132" = 71457.69
  "L132" = 29739.858
  "M132" = -68927.69
  "N132" = -34799.858
  "H136" = 4073.6667
  "I136" = 3832.1458
  "K136" = 11496.4374
  "M136" = -8946.437399999999
}
foreach ($ref in $ARM_values.Keys) {
    $ws.Range($ref).Value = $ARM_values[$ref]
}

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$BSM_values = @{
  "H3" = 898.84
  "I3" = 821.5454999999999
  "K3" = 821.5454999999999
  "M3" = -707.5454999999999
  "H59" = 99998
  "J59" = 99998
  "L59" = 99998
  "N59" = -101692
  "H105" = 2052.5625
  "I105" = 2052.5625
  "K105" = 2052.5625
  "M105" = -305.5625
  "H134" = 3271
  "I134" = 3656.9
  "J134" = 2995.3572
  "K134" = 10970.7
  "L134" = 8986.071599999999
  "M134" = -8435.700000000001
  "N134" = -14056.0716
}
foreach ($ref in $BSM_values.Keys) {
    $ws.Range($ref).Value = $BSM_values[$ref]
}

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$CRP_values = @{
  "H16" = 1608.4
  "I16" = 1608.4
  "K16" = 1608.4
  "M16" = -1321.4
  "H97" = 0
  "J97" = 0
  "L97" = 0
  "H113" = 1608.4
  "I113" = 1608.4
  "K113" = 1608.4
  "M113" = 561.5999999999999
  "H132" = 40405760
  "I132" = 45978610
  "K132" = 137935830
  "M132" = -137933300
  "H134" = 3330.4375
  "I134" = 2750.3
  "K134" = 8250.900000000001
  "M134" = -5715.900000000001
  "H141" = 164342.44
  "J141" = 194197.89
  "L141" = 194197.89
  "N141" = -204557.89
}
foreach ($ref in $CRP_values.Keys) {
    $ws.Range($ref).Value = $CRP_values[$ref]
}
$CRP_deletes = @("N97")
foreach ($ref in $CRP_deletes) {
    $ws.Range($ref).ClearContents()
}

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$CUL_values = @{
  "H4" = 5271966
  "I4" = 347190.62
  "K4" = 1041571.86
  "M4" = -1041459.86
  "H37" = 1000000000
  "J37" = 1000000000
  "L37" = 3000000000
  "N37" = -3000000224
  "H38" = 263.89474
  "J38" = 435.55554
  "L38" = 1306.66662
  "N38" = -2000.66662
}
foreach ($ref in $CUL_values.Keys) {
    $ws.Range($ref).Value = $CUL_values[$ref]
}

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$GSM_values = @{
  "H2" = 365.36
  "I2" = 173.07692
  "K2" = 173.07692
  "M2" = -60.07692
  "H102" = 442015.34
  "I102" = 516966.9
  "J102" = 3013.4285
  "K102" = 516966.9
  "L102" = 3013.4285
  "M102" = -515344.9
  "N102" = -6257.4285
}
foreach ($ref in $GSM_values.Keys) {
    $ws.Range($ref).Value = $GSM_values[$ref]
}

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$LTW_values = @{
  "H7" = 7238.4443
  "I7" = 6691.1665
  "J7" = 8333
  "K7" = 6691.1665
  "L7" = 8333
  "M7" = -6579.1665
  "N7" = -8557
  "H40" = 8931588
  "I40" = 10207131
  "K40" = 10207131
  "M40" = -10206995
  "H46" = 4646.4736
  "I46" = 4772.5
  "J46" = 4612.8667
  "K46" = 4772.5
  "L46" = 4612.8667
  "M46" = -4584.5
  "N46" = -4988.8667
  "H126" = 7238.4443
  "I126" = 6691.1665
  "J126" = 8333
  "K126" = 20073.4995
  "L126" = 24999
  "M126" = -17603.4995
  "N126" = -29939
}
foreach ($ref in $LTW_values.Keys) {
    $ws.Range($ref).Value = $LTW_values[$ref]
}

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$WVR_values = @{
  "H113" = 1678.4286
  "I113" = 1433
  "J113" = 1862.5
  "K113" = 4299
  "L113" = 5587.5
  "M113" = -2129
  "N113" = -9927.5
  "H126" = 47621310
  "I126" = 58825520
  "K126" = 176476560
  "M126" = -176474090
  "H132" = 2096.2368
  "I132" = 1034.6086
  "K132" = 3103.8258
  "M132" = -573.8258000000001
  "H136" = 3091.4524
  "I136" = 2159.1875
  "J136" = 6074.7
  "K136" = 6477.5625
  "L136" = 18224.1
  "M136" = -3927.5625
  "N136" = -23324.1
  "H138" = 78943
}
foreach ($ref in $WVR_values.Keys) {
    $ws.Range($ref).Value = $WVR_values[$ref]
}
